$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 9: Sabertooth/Syren entry
$ws.Range("A9").Value = "Sabertooth/Syren "
$ws.Range("B9").Value = $ws.Range("B2").Value2
$ws.Range("E9").Value = "http://www.dimensionengineering.com/info/arduino"
$ws.Range("F9").Value = "Drivers for Syren 10 and Sabertooth 2x25"

# Copy the date formatting from an existing "Last Downloaded" cell so the
# new cell reuses the same style (numFmtId 14, mm/dd/yyyy) already in the
# workbook rather than creating a new one.
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C9").Value = (Get-Date -Year 2014 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0)

# Column F (Description) needs to grow to fit the new, longer description text
$ws.Columns.Item(6).ColumnWidth = 36.6

# Move active selection to G2, matching the recorded cursor position after edit
$ws.Range("G2").Select()
